# Docx writer: Use different style for block quotes in notes.
# Add a new "Footnote Block Text" paragraph style, based on "Footnote Text",
# mirroring how the existing "Block Text" style is derived from "Body Text".

$d = $word.ActiveDocument

# 1 = wdStyleTypeParagraph
$s = $d.Styles.Add("FootnoteBlockText", 1)
$s.NameLocal = "Footnote Block Text"
$s.BaseStyle = "Footnote Text"
$s.NextParagraphStyle = "Footnote Text"
$s.Priority = 9
$s.UnhideWhenUsed = $true
$s.QuickStyle = $true

# ParagraphFormat distances are expressed in points; the target spacing/indents
# are 100 and 480 twentieths-of-a-point (i.e. 5pt and 24pt), same as BlockText.
$s.ParagraphFormat.SpaceBefore = 5
$s.ParagraphFormat.SpaceAfter = 5
$s.ParagraphFormat.FirstLineIndent = 0
$s.ParagraphFormat.LeftIndent = 24
$s.ParagraphFormat.RightIndent = 24
